$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 3-5 (cyclic shift of F:V content) and rows 70-71 (swap of F:V content) ---
# Row 3
$ws.Range("F3").Value = "Orijent"
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = "Dubrava"
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 2.12
$ws.Range("K3").Value = "11/08/2023 05:43"
$ws.Range("L3").Value = 2.01
$ws.Range("M3").Value = "12/08/2023 17:29"
$ws.Range("N3").Value = 3.2
$ws.Range("O3").Value = "11/08/2023 05:43"
$ws.Range("P3").Value = 3.38
$ws.Range("Q3").Value = "12/08/2023 17:29"
$ws.Range("R3").Value = 3.12
$ws.Range("S3").Value = "11/08/2023 05:43"
$ws.Range("T3").Value = 3.59
$ws.Range("U3").Value = "12/08/2023 17:29"
$ws.Range("V3").Value = "https://www.betexplorer.com/football/croatia/prva-nl/orijent-dubrava-zagreb/Cdce1OPO/"

# Row 4
$ws.Range("F4").Value = "Cibalia"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = "Croatia Zmijavci"
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 1.65
$ws.Range("K4").Value = "11/08/2023 05:43"
$ws.Range("L4").Value = 2.02
$ws.Range("M4").Value = "12/08/2023 15:37"
$ws.Range("N4").Value = 3.47
$ws.Range("O4").Value = "11/08/2023 05:43"
$ws.Range("P4").Value = 3.47
$ws.Range("Q4").Value = "12/08/2023 17:16"
$ws.Range("R4").Value = 4.52
$ws.Range("S4").Value = "11/08/2023 05:43"
$ws.Range("T4").Value = 3.45
$ws.Range("U4").Value = "12/08/2023 17:16"
$ws.Range("V4").Value = "https://www.betexplorer.com/football/croatia/prva-nl/cibalia-croatia-zmijavci/E1yvQ72P/"

# Row 5
$ws.Range("F5").Value = "Sesvete"
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = "Dugopolje"
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 2.49
$ws.Range("K5").Value = "12/08/2023 17:11"
$ws.Range("L5").Value = 2.26
$ws.Range("M5").Value = "12/08/2023 17:25"
$ws.Range("N5").Value = 2.94
$ws.Range("O5").Value = "12/08/2023 17:11"
$ws.Range("P5").Value = 3.14
$ws.Range("Q5").Value = "12/08/2023 17:25"
$ws.Range("R5").Value = 2.65
$ws.Range("S5").Value = "12/08/2023 17:11"
$ws.Range("T5").Value = 3.23
$ws.Range("U5").Value = "12/08/2023 17:25"
$ws.Range("V5").Value = "https://www.betexplorer.com/football/croatia/prva-nl/sesvete-dugopolje/jNS5LRAt/"

# Row 70
$ws.Range("F70").Value = "Bijelo Brdo"
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = "Solin"
$ws.Range("I70").Value = 1
$ws.Range("J70").Value = 2.14
$ws.Range("K70").Value = "27/10/2023 03:12"
$ws.Range("L70").Value = 2.68
$ws.Range("M70").Value = "28/10/2023 14:57"
$ws.Range("N70").Value = 3.17
$ws.Range("O70").Value = "27/10/2023 03:12"
$ws.Range("P70").Value = 2.66
$ws.Range("Q70").Value = "28/10/2023 14:56"
$ws.Range("R70").Value = 3.02
$ws.Range("S70").Value = "27/10/2023 03:12"
$ws.Range("T70").Value = 3.1
$ws.Range("U70").Value = "28/10/2023 14:57"
$ws.Range("V70").Value = "https://www.betexplorer.com/football/croatia/prva-nl/bijelo-brdo-solin/hM3tS0Qd/"

# Row 71
$ws.Range("F71").Value = "Dugopolje"
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = "Sesvete"
$ws.Range("I71").Value = 1
$ws.Range("J71").Value = 1.72
$ws.Range("K71").Value = "27/10/2023 03:12"
$ws.Range("L71").Value = 1.68
$ws.Range("M71").Value = "27/10/2023 13:14"
$ws.Range("N71").Value = 3.47
$ws.Range("O71").Value = "27/10/2023 03:12"
$ws.Range("P71").Value = 3.76
$ws.Range("Q71").Value = "28/10/2023 14:26"
$ws.Range("R71").Value = 4.07
$ws.Range("S71").Value = "27/10/2023 03:12"
$ws.Range("T71").Value = 4.74
$ws.Range("U71").Value = "28/10/2023 14:26"
$ws.Range("V71").Value = "https://www.betexplorer.com/football/croatia/prva-nl/dugopolje-sesvete/Yya1NxuS/"

# --- Append new rows 81-83 ---
# Copy formatting (styles) of row 80 (A:V) down into the three new rows first
$ws.Range("A80:V80").Copy()
$ws.Range("A81:V81").PasteSpecial(-4122)
$ws.Range("A82:V82").PasteSpecial(-4122)
$ws.Range("A83:V83").PasteSpecial(-4122)

# Row 81
$ws.Range("A81").Value = 80
$ws.Range("B81").Value = "croatia"
$ws.Range("C81").Value = "prva-nl"
$ws.Range("D81").Value = "2023-2024"
$ws.Range("E81").Value = 45241.57638888889
$ws.Range("F81").Value = "Dubrava"
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = "Sesvete"
$ws.Range("I81").Value = 2
$ws.Range("J81").Value = 1.95
$ws.Range("K81").Value = "10/11/2023 02:12"
$ws.Range("L81").Value = 2.1
$ws.Range("M81").Value = "11/11/2023 13:46"
$ws.Range("N81").Value = 3.33
$ws.Range("O81").Value = "10/11/2023 02:12"
$ws.Range("P81").Value = 3.3
$ws.Range("Q81").Value = "11/11/2023 13:43"
$ws.Range("R81").Value = 3.3
$ws.Range("S81").Value = "10/11/2023 02:12"
$ws.Range("T81").Value = 3.44
$ws.Range("U81").Value = "11/11/2023 13:46"
$ws.Range("V81").Value = "https://www.betexplorer.com/football/croatia/prva-nl/dubrava-zagreb-sesvete/h4AAVZWl/"

# Row 82
$ws.Range("A82").Value = 81
$ws.Range("B82").Value = "croatia"
$ws.Range("C82").Value = "prva-nl"
$ws.Range("D82").Value = "2023-2024"
$ws.Range("E82").Value = 45241.58333333334
$ws.Range("F82").Value = "Bijelo Brdo"
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = "Sibenik"
$ws.Range("I82").Value = 1
$ws.Range("J82").Value = 4.54
$ws.Range("K82").Value = "10/11/2023 02:12"
$ws.Range("L82").Value = 5.86
$ws.Range("M82").Value = "11/11/2023 13:52"
$ws.Range("N82").Value = 3.58
$ws.Range("O82").Value = "10/11/2023 02:12"
$ws.Range("P82").Value = 4.01
$ws.Range("Q82").Value = "11/11/2023 13:52"
$ws.Range("R82").Value = 1.63
$ws.Range("S82").Value = "10/11/2023 02:12"
$ws.Range("T82").Value = 1.54
$ws.Range("U82").Value = "11/11/2023 13:51"
$ws.Range("V82").Value = "https://www.betexplorer.com/football/croatia/prva-nl/bijelo-brdo-sibenik/z5LWxIvF/"

# Row 83
$ws.Range("A83").Value = 82
$ws.Range("B83").Value = "croatia"
$ws.Range("C83").Value = "prva-nl"
$ws.Range("D83").Value = "2023-2024"
$ws.Range("E83").Value = 45241.58333333334
$ws.Range("F83").Value = "Vukovar 1991"
$ws.Range("G83").Value = 2
$ws.Range("H83").Value = "Orijent"
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1.47
$ws.Range("K83").Value = "10/11/2023 02:12"
$ws.Range("L83").Value = 1.53
$ws.Range("M83").Value = "11/11/2023 13:54"
$ws.Range("N83").Value = 4.22
$ws.Range("O83").Value = "10/11/2023 02:12"
$ws.Range("P83").Value = 4.56
$ws.Range("Q83").Value = "11/11/2023 13:54"
$ws.Range("R83").Value = 5.03
$ws.Range("S83").Value = "10/11/2023 02:12"
$ws.Range("T83").Value = 5
$ws.Range("U83").Value = "11/11/2023 13:54"
$ws.Range("V83").Value = "https://www.betexplorer.com/football/croatia/prva-nl/vukovar-1991-orijent/8WYwyd9R/"

$excel.CutCopyMode = $false

